$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 19-55: each entry is (row, terminal-group name, charger name, last-charge-end serial datetime)
$rows = @(
    @(19, "长沙特来电飞狐四方坪西区充电站", "604号直流", 45939.747453703705),
    @(20, "长沙特来电飞狐四方坪东区充电站", "006B号直流", 45940.517048611109),
    @(21, "长沙特来电飞狐四方坪东区充电站", "003B号直流", 45940.552430555559),
    @(22, "长沙特来电飞狐四方坪西区充电站", "605号直流", 45940.556215277778),
    @(23, "长沙特来电飞狐四方坪南区充电站", "406号直流", 45940.557627314818),
    @(24, "长沙特来电飞狐四方坪西区充电站", "603号直流", 45940.558483796296),
    @(25, "长沙特来电飞狐四方坪西区充电站", "602号直流", 45941.043564814812),
    @(26, "长沙特来电飞狐四方坪西区充电站", "B04号直流", 45941.079953703702),
    @(27, "长沙特来电飞狐四方坪西区充电站", "405号直流", 45941.19804398148),
    @(28, "长沙特来电飞狐四方坪西区充电站", "403号直流", 45941.233275462961),
    @(29, "长沙特来电飞狐四方坪南区充电站", "103号直流", 45941.259837962964),
    @(30, "长沙特来电飞狐四方坪西区充电站", "404号直流", 45941.277685185189),
    @(31, "长沙市开福区高岭香江国际城充电站建设项目", "108号直流", 45941.331342592595),
    @(32, "长沙市开福区高岭香江国际城充电站建设项目", "110号直流", 45941.345983796295),
    @(33, "长沙特来电飞狐四方坪西区充电站", "901号直流", 45941.365798611114),
    @(34, "长沙特来电飞狐四方坪东区充电站", "001A号直流", 45941.527708333335),
    @(35, "长沙特来电飞狐四方坪东区充电站", "102号直流", 45941.538611111115),
    @(36, "长沙市开福区高岭香江国际城充电站建设项目", "305号直流", 45941.542245370372),
    @(37, "长沙特来电飞狐四方坪东区充电站", "402号直流", 45941.545925925922),
    @(38, "长沙市开福区高岭香江国际城充电站建设项目", "203号直流", 45941.555543981478),
    @(39, "长沙特来电飞狐四方坪南区充电站", "305号直流", 45941.55740740741),
    @(40, "长沙特来电飞狐四方坪南区充电站", "201号直流", 45941.561180555553),
    @(41, "长沙特来电飞狐四方坪西区充电站", "401号直流", 45941.584722222222),
    @(42, "长沙特来电飞狐四方坪西区充电站", "303号直流", 45941.591643518521),
    @(43, "长沙特来电飞狐四方坪东区充电站", "401号直流", 45941.595682870371),
    @(44, "长沙市开福区高岭香江国际城充电站建设项目", "303号直流", 45941.641412037039),
    @(45, "长沙市开福区高岭香江国际城充电站建设项目", "105号直流", 45941.650648148148),
    @(46, "长沙市开福区高岭香江国际城充电站建设项目", "106号直流", 45941.65084490741),
    @(47, "长沙市开福区高岭香江国际城充电站建设项目", "101号直流", 45941.66269675926),
    @(48, "长沙市开福区高岭香江国际城充电站建设项目", "103号直流", 45941.666030092594),
    @(49, "长沙特来电飞狐四方坪西区充电站", "905号直流", 45941.682395833333),
    @(50, "长沙特来电飞狐四方坪南区充电站", "403号直流", 45941.728981481479),
    @(51, "长沙市开福区高岭香江国际城充电站建设项目", "208号直流", 45941.730914351851),
    @(52, "长沙特来电飞狐四方坪南区充电站", "202号直流", 45941.732175925928),
    @(53, "长沙特来电飞狐四方坪东区充电站", "404号直流", 45941.73809027778),
    @(54, "长沙市开福区高岭香江国际城充电站建设项目", "211号直流", 45941.752488425926),
    @(55, "长沙市开福区高岭香江国际城充电站建设项目", "112号直流", 45941.755798611113)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

$ws.Range("D16").Select()
